$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 9052.6
$ws.Range("I82").Value = 7566.25
$ws.Range("K82").Value = 22698.75
$ws.Range("M82").Value = -22292.75
$ws.Range("H85").Value = 9052.6
$ws.Range("I85").Value = 7566.25
$ws.Range("K85").Value = 22698.75
$ws.Range("M85").Value = -21294.75
$ws.Range("H96").Value = 111120830
$ws.Range("I96").Value = 7499.3335
$ws.Range("J96").Value = 166677500
$ws.Range("K96").Value = 22498.0005
$ws.Range("L96").Value = 500032500
$ws.Range("M96").Value = -21125.0005
$ws.Range("N96").Value = -500035246
$ws.Range("H111").Value = 952.25
$ws.Range("I111").Value = 918.0769
$ws.Range("J111").Value = 1100.3334
$ws.Range("K111").Value = 2754.2307
$ws.Range("L111").Value = 3301.0002
$ws.Range("M111").Value = 312.7692999999999
$ws.Range("N111").Value = -9435.0002
$ws.Range("H112").Value = 4418.5713
$ws.Range("J112").Value = 3523.6365
$ws.Range("L112").Value = 10570.9095
$ws.Range("N112").Value = -12786.9095
$ws.Range("H118").Value = 648
$ws.Range("I118").Value = 295.85715
$ws.Range("K118").Value = 887.5714499999999
$ws.Range("M118").Value = 769.4285500000001
$ws.Range("I129").Value = 797.3333
$ws.Range("K129").Value = 2391.9999
$ws.Range("M129").Value = 2608.0001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 16712689
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H32").Value = 9363.34
$ws.Range("I32").Value = 4701.184
$ws.Range("J32").Value = 24126.834
$ws.Range("K32").Value = 4701.184
$ws.Range("L32").Value = 24126.834
$ws.Range("M32").Value = -4414.184
$ws.Range("N32").Value = -24700.834
$ws.Range("H61").Value = 9727.75
$ws.Range("I61").Value = 11353.667
$ws.Range("K61").Value = 11353.667
$ws.Range("M61").Value = -11141.667
$ws.Range("H74").Value = 5428.615
$ws.Range("I74").Value = 5553.5
$ws.Range("K74").Value = 5553.5
$ws.Range("M74").Value = -4679.5
$ws.Range("H77").Value = 5428.615
$ws.Range("I77").Value = 5553.5
$ws.Range("K77").Value = 27767.5
$ws.Range("M77").Value = -23399.5
$ws.Range("H98").Value = 37500
$ws.Range("J98").Value = 37500
$ws.Range("L98").Value = 37500
$ws.Range("N98").Value = -43490
$ws.Range("H122").Value = 6543.8286
$ws.Range("I122").Value = 5127.6523
$ws.Range("J122").Value = 9258.166999999999
$ws.Range("K122").Value = 15382.9569
$ws.Range("L122").Value = 27774.501
$ws.Range("M122").Value = -12932.9569
$ws.Range("N122").Value = -32674.501
$ws.Range("H136").Value = 9727.75
$ws.Range("I136").Value = 11353.667
$ws.Range("K136").Value = 34061.001
$ws.Range("M136").Value = -31511.001
$ws.Range("N13").ClearContents()

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 19597.223
$ws.Range("I134").Value = 11319.4
$ws.Range("J134").Value = 29944.5
$ws.Range("K134").Value = 33958.2
$ws.Range("L134").Value = 89833.5
$ws.Range("M134").Value = -31423.2
$ws.Range("N134").Value = -94903.5
$ws.Range("H141").Value = 171501.97
$ws.Range("J141").Value = 171501.97
$ws.Range("L141").Value = 171501.97
$ws.Range("N141").Value = -181861.97

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2126.4546
$ws.Range("I34").Value = 297.2143
$ws.Range("J34").Value = 5327.625
$ws.Range("K34").Value = 891.6428999999999
$ws.Range("L34").Value = 15982.875
$ws.Range("M34").Value = -807.6428999999999
$ws.Range("N34").Value = -16150.875
$ws.Range("H94").Value = 13999.333
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("H99").Value = 4583
$ws.Range("I99").Value = 3499.8
$ws.Range("J99").Value = 9999
$ws.Range("K99").Value = 10499.4
$ws.Range("L99").Value = 29997
$ws.Range("M99").Value = -8253.400000000001
$ws.Range("N99").Value = -34489
$ws.Range("H124").Value = 16249.6
$ws.Range("I124").Value = 7812
$ws.Range("J124").Value = 50000
$ws.Range("K124").Value = 23436
$ws.Range("L124").Value = 150000
$ws.Range("M124").Value = -18526
$ws.Range("N124").Value = -159820
$ws.Range("H139").Value = 33336192
$ws.Range("I139").Value = 50002000
$ws.Range("J139").Value = 4578.4
$ws.Range("K139").Value = 150006000
$ws.Range("L139").Value = 13735.2
$ws.Range("M139").Value = -150000860
$ws.Range("N139").Value = -24015.2
$ws.Range("H140").Value = 34093624
$ws.Range("I140").Value = 39476196
$ws.Range("K140").Value = 118428588
$ws.Range("M140").Value = -118423408
$ws.Range("M94").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H123").Value = 61600
$ws.Range("J123").Value = 61600
$ws.Range("L123").Value = 61600
$ws.Range("N123").Value = -66500
$ws.Range("H126").Value = 46166404
$ws.Range("I126").Value = 166668140
$ws.Range("K126").Value = 500004420
$ws.Range("M126").Value = -500001950
$ws.Range("H132").Value = 11607.7
$ws.Range("I132").Value = 12203
$ws.Range("K132").Value = 36609
$ws.Range("M132").Value = -34079
$ws.Range("H135").Value = 109722
$ws.Range("J135").Value = 109722
$ws.Range("L135").Value = 109722
$ws.Range("N135").Value = -119862
$ws.Range("N96").ClearContents()

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6421.75
$ws.Range("I40").Value = 7397.6665
$ws.Range("J40").Value = 5445.8335
$ws.Range("K40").Value = 7397.6665
$ws.Range("L40").Value = 5445.8335
$ws.Range("M40").Value = -7261.6665
$ws.Range("N40").Value = -5717.8335
$ws.Range("H46").Value = 41668056
$ws.Range("I46").Value = 1009
$ws.Range("J46").Value = 71430230
$ws.Range("K46").Value = 1009
$ws.Range("L46").Value = 71430230
$ws.Range("M46").Value = -821
$ws.Range("N46").Value = -71430606
$ws.Range("H122").Value = 4645.9062
$ws.Range("I122").Value = 4049.5454
$ws.Range("J122").Value = 5957.9
$ws.Range("K122").Value = 12148.6362
$ws.Range("L122").Value = 17873.7
$ws.Range("M122").Value = -9698.636200000001
$ws.Range("N122").Value = -22773.7
$ws.Range("H132").Value = 5651.143
$ws.Range("I132").Value = 4097.6895
$ws.Range("K132").Value = 12293.0685
$ws.Range("M132").Value = -9763.068500000001
$ws.Range("H136").Value = 66683320
$ws.Range("I136").Value = 12483.333
$ws.Range("K136").Value = 37449.999
$ws.Range("M136").Value = -34899.999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1857.3334
$ws.Range("J81").Value = 4624.5
$ws.Range("L81").Value = 9249
$ws.Range("N81").Value = -11371
$ws.Range("H84").Value = 1857.3334
$ws.Range("J84").Value = 4624.5
$ws.Range("L84").Value = 46245
$ws.Range("N84").Value = -56853
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H122").Value = 12516.5
$ws.Range("I122").Value = 3904.2
$ws.Range("K122").Value = 11712.6
$ws.Range("M122").Value = -9262.599999999999
$ws.Range("H126").Value = 6755.231
$ws.Range("I126").Value = 4366.857
$ws.Range("K126").Value = 13100.571
$ws.Range("M126").Value = -10630.571
$ws.Range("H132").Value = 11070.857
$ws.Range("I132").Value = 10684.363
$ws.Range("K132").Value = 32053.089
$ws.Range("M132").Value = -29523.089
$ws.Range("H136").Value = 25021032
$ws.Range("J136").Value = 17429.75
$ws.Range("L136").Value = 52289.25
$ws.Range("N136").Value = -57389.25
$ws.Range("N92").ClearContents()
